# Updates cryptos list with fresh price/volume data (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.231.03"
$ws.Range("D3").Value = "3.117.69"
$ws.Range("E3").Value = "  +3.57%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.18%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.110.76"
$ws.Range("E8").Value = "  +3.76%  "
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.62%  "
$ws.Range("E12").Value = "  +2.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.56%  "
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "3.632.87"
$ws.Range("E16").Value = "  +3.59%  "
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "63.126.42"
$ws.Range("E18").Value = "  +5.70%  "
$ws.Range("D19").Value = "3.113.55"
$ws.Range("E19").Value = "  +3.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.14%  "
$ws.Range("E21").Value = "  +3.08%  "
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("E23").Value = "  +5.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +4.84%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.64%  "
$ws.Range("E31").Value = "  +8.48%  "
$ws.Range("E32").Value = "  +4.15%  "
$ws.Range("E33").Value = "  +4.99%  "
$ws.Range("D34").Value = "0.0₃0871"
$ws.Range("E34").Value = "  +10.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.86%  "
$ws.Range("E36").Value = "  +4.67%  "
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +16.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "433.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.49%  "
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("D42").Value = "2.922.58"
$ws.Range("E42").Value = "  +5.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0370"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.28%  "
$ws.Range("E44").Value = "  +10.00%  "
$ws.Range("E45").Value = "  +5.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.97%  "
